$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.01135533333333333
$ws.Range("H2").Value = 0.034066
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.2615913333333333
$ws.Range("N2").Value = 0.784774
$ws.Range("O2").Value = 0.08239613548481725
$ws.Range("P2").Value = 0.08239613548481727
$ws.Range("Q2").Value = 0.002970456787111111
$ws.Range("R2").Value = 0.026734111084
$ws.Range("S2").Value = 0.08239613548481725
$ws.Range("T2").Value = 0.08239613548481727

# Row 3
$ws.Range("G3").Value = 0.01135533333333333
$ws.Range("H3").Value = 0.034066
$ws.Range("N3").Value = 5.233242000000001
$ws.Range("O3").Value = 0.5494561706387266
$ws.Range("P3").Value = 0.5494561706387268
$ws.Range("Q3").Value = 0.01980840244133333
$ws.Range("R3").Value = 0.178275621972
$ws.Range("S3").Value = 0.5494561706387266
$ws.Range("T3").Value = 0.5494561706387268

# Row 4
$ws.Range("G4").Value = 0.01135533333333333
$ws.Range("H4").Value = 0.034066
$ws.Range("M4").Value = 1.168795666666667
$ws.Range("N4").Value = 3.506387
$ws.Range("O4").Value = 0.3681476938764561
$ws.Range("P4").Value = 0.3681476938764561
$ws.Range("Q4").Value = 0.01327206439355556
$ws.Range("R4").Value = 0.119448579542
$ws.Range("S4").Value = 0.3681476938764561
$ws.Range("T4").Value = 0.3681476938764561
